# Gant - Ingeniería WEB.xlsx edits
#  - Mark columns G:I as "Realizado" (done) for rows 10-12 (previously "Pendiente")
#  - Move the "today" connector line marker from column G to column J
#  - Update the active selection to I12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update progress cells: G10:I10, G11:I11, G12:I12 -> "Realizado" / "Bueno" style ---
$doneRanges = @("G10:I10", "G11:I11", "G12:I12")
foreach ($r in $doneRanges) {
    $rng = $ws.Range($r)
    $rng.Value2 = "Realizado"
    $rng.Style = "Bueno"
}

# --- Move the connector (today-marker) line shape to its new position ---
$line = $ws.Shapes.Item("Conector recto 2")
$line.Left = 694.5
$line.Top = 57.75
$line.Width = 0
$line.Height = 206.25

# --- Update the active selection shown when the workbook is opened ---
$ws.Range("I12").Select()
